# Daily attendance processing - 2025-12-17 07:08:07
# Rotate the "Recorded By" (column G) comma-separated list left by one
# position (move the first name/email to the end) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -le 1) { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "

    if ($rotated -ne $text) {
        $cell.Value2 = $rotated
    }
}
